# Scheduled-runner update: refresh cached market-board price/profit figures
# across all job sheets. Values below are literal (no formulas in source).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H6").Value  = 734.7778
$ws.Range("I6").Value  = 512.1667
$ws.Range("J6").Value  = 1180
$ws.Range("K6").Value  = 1536.5001
$ws.Range("L6").Value  = 3540
$ws.Range("M6").Value  = -1424.5001
$ws.Range("N6").Value  = -3764

$ws.Range("H17").Value = 1000385.06
$ws.Range("J17").Value = 1000385.06
$ws.Range("L17").Value = 3001155.18
$ws.Range("N17").Value = -3001491.18

$ws.Range("H74").Value = 3811.2778
$ws.Range("I74").Value = 3550.25
$ws.Range("J74").Value = 4333.3335
$ws.Range("K74").Value = 3550.25
$ws.Range("L74").Value = 4333.3335
$ws.Range("M74").Value = -2614.25
$ws.Range("N74").Value = -6205.3335

$ws.Range("H77").Value = 3811.2778
$ws.Range("I77").Value = 3550.25
$ws.Range("J77").Value = 4333.3335
$ws.Range("K77").Value = 17751.25
$ws.Range("L77").Value = 21666.6675
$ws.Range("M77").Value = -13071.25
$ws.Range("N77").Value = -31026.6675

$ws.Range("H125").Value = 736700.3
$ws.Range("I125").Value = 1651.5555
$ws.Range("J125").Value = 1681763
$ws.Range("K125").Value = 14863.9995
$ws.Range("L125").Value = 15135867
$ws.Range("M125").Value = -12403.9995
$ws.Range("N125").Value = -15140787

$ws.Range("H129").Value = 1092.7755
$ws.Range("J129").Value = 1092.7755
$ws.Range("L129").Value = 3278.3265
$ws.Range("N129").Value = -13278.3265

$ws.Range("H137").Value = 973.64
$ws.Range("I137").Value = 889.975
$ws.Range("J137").Value = 1308.3
$ws.Range("K137").Value = 2669.925
$ws.Range("L137").Value = 3924.9
$ws.Range("M137").Value = -119.9250000000002
$ws.Range("N137").Value = -9024.9

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 50000
$ws.Range("K3").Value = 50000
$ws.Range("M3").Value = -49885

$ws.Range("H102").Value = 1282.9
$ws.Range("I102").Value = 1141.125
$ws.Range("J102").Value = 1850
$ws.Range("K102").Value = 1141.125
$ws.Range("L102").Value = 1850
$ws.Range("M102").Value = 480.875
$ws.Range("N102").Value = -5094

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value = 10303.23
$ws.Range("I20").Value = 3203.2
$ws.Range("J20").Value = 14740.75
$ws.Range("K20").Value = 3203.2
$ws.Range("L20").Value = 14740.75
$ws.Range("M20").Value = -2956.2
$ws.Range("N20").Value = -15234.75

$ws.Range("H105").Value = 1978.5385
$ws.Range("I105").Value = 2200.9092
$ws.Range("J105").Value = 755.5
$ws.Range("K105").Value = 2200.9092
$ws.Range("L105").Value = 755.5
$ws.Range("M105").Value = -453.9092000000001
$ws.Range("N105").Value = -4249.5

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H3").Value = 20334.25
$ws.Range("J3").Value = 53335.332
$ws.Range("L3").Value = 53335.332
$ws.Range("N3").Value = -53561.332

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -580

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H5").Value = 12801
$ws.Range("J5").Value = 12801
$ws.Range("L5").Value = 12801
$ws.Range("N5").Value = -13025

$ws.Range("H80").Value = 4700
$ws.Range("I80").Value = 5600
$ws.Range("J80").Value = 3800
$ws.Range("K80").Value = 5600
$ws.Range("L80").Value = 3800
$ws.Range("M80").Value = -4602
$ws.Range("N80").Value = -5796

$ws.Range("H83").Value = 4700
$ws.Range("I83").Value = 5600
$ws.Range("J83").Value = 3800
$ws.Range("K83").Value = 28000
$ws.Range("L83").Value = 19000
$ws.Range("M83").Value = -23008
$ws.Range("N83").Value = -28984

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 93327.55
$ws.Range("I7").Value = 112844.78
$ws.Range("J7").Value = 5500
$ws.Range("K7").Value = 112844.78
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = -112732.78
$ws.Range("N7").Value = -5724

$ws.Range("H22").Value = 569.9286
$ws.Range("I22").Value = 581.2857
$ws.Range("J22").Value = 558.5714
$ws.Range("K22").Value = 581.2857
$ws.Range("L22").Value = 558.5714
$ws.Range("M22").Value = -286.2857
$ws.Range("N22").Value = -1148.5714

$ws.Range("H27").Value = 569.9286
$ws.Range("I27").Value = 581.2857
$ws.Range("J27").Value = 558.5714
$ws.Range("K27").Value = 581.2857
$ws.Range("L27").Value = 558.5714
$ws.Range("M27").Value = -474.2857
$ws.Range("N27").Value = -772.5714

$ws.Range("H40").Value = 25913.5
$ws.Range("I40").Value = 27554.85
$ws.Range("J40").Value = 9500
$ws.Range("K40").Value = 27554.85
$ws.Range("L40").Value = 9500
$ws.Range("M40").Value = -27418.85
$ws.Range("N40").Value = -9772

$ws.Range("H46").Value = 111921.336
$ws.Range("I46").Value = 200558.4
$ws.Range("J46").Value = 1125
$ws.Range("K46").Value = 200558.4
$ws.Range("L46").Value = 1125
$ws.Range("M46").Value = -200370.4
$ws.Range("N46").Value = -1501

$ws.Range("H68").Value = 1482.5
$ws.Range("I68").Value = 1233.3334
$ws.Range("J68").Value = 1632
$ws.Range("K68").Value = 1233.3334
$ws.Range("L68").Value = 1632
$ws.Range("M68").Value = -484.3334
$ws.Range("N68").Value = -3130

$ws.Range("H71").Value = 1482.5
$ws.Range("I71").Value = 1233.3334
$ws.Range("J71").Value = 1632
$ws.Range("K71").Value = 6166.666999999999
$ws.Range("L71").Value = 8160
$ws.Range("M71").Value = -2422.666999999999
$ws.Range("N71").Value = -15648

$ws.Range("H82").Value = 3089.3
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 3611.625
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 3611.625
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -4333.625

$ws.Range("H85").Value = 3089.3
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 3611.625
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 3611.625
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -6107.625

$ws.Range("H122").Value = 13890601
$ws.Range("I122").Value = 27778832
$ws.Range("J122").Value = 2370
$ws.Range("K122").Value = 83336496
$ws.Range("L122").Value = 7110
$ws.Range("M122").Value = -83334046
$ws.Range("N122").Value = -12010

$ws.Range("H126").Value = 93327.55
$ws.Range("I126").Value = 112844.78
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 338534.34
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -336064.34
$ws.Range("N126").Value = -21440

$ws.Range("H132").Value = 5345.7866
$ws.Range("I132").Value = 5134.5093
$ws.Range("K132").Value = 15403.5279
$ws.Range("M132").Value = -12873.5279

$ws.Range("H135").Value = 66174.14
$ws.Range("J135").Value = 66174.14
$ws.Range("L135").Value = 66174.14
$ws.Range("N135").Value = -76314.14

$ws.Range("H136").Value = 12822611
$ws.Range("I136").Value = 2182
$ws.Range("J136").Value = 333333340
$ws.Range("K136").Value = 6546
$ws.Range("L136").Value = 1000000020
$ws.Range("M136").Value = -3996
$ws.Range("N136").Value = -1000005120

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 5712.857
$ws.Range("J62").Value = 5712.857
$ws.Range("L62").Value = 5712.857
$ws.Range("N62").Value = -6960.857

$ws.Range("H65").Value = 5712.857
$ws.Range("J65").Value = 5712.857
$ws.Range("L65").Value = 28564.285
$ws.Range("N65").Value = -34804.285

$ws.Range("H81").Value = 1926.1818
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 1968.8
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 3937.6
$ws.Range("M81").Value = -1939
$ws.Range("N81").Value = -6059.6

$ws.Range("H84").Value = 1926.1818
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 1968.8
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 19688
$ws.Range("M84").Value = -9696
$ws.Range("N84").Value = -30296

$ws.Range("H96").Value = 1179.25
$ws.Range("I96").Value = 866.6667
$ws.Range("J96").Value = 1366.8
$ws.Range("K96").Value = 866.6667
$ws.Range("L96").Value = 1366.8
$ws.Range("M96").Value = 506.3333
$ws.Range("N96").Value = -4112.8

$ws.Range("H122").Value = 43488.957
$ws.Range("I122").Value = 60425.94
$ws.Range("J122").Value = 2356.2856
$ws.Range("K122").Value = 181277.82
$ws.Range("L122").Value = 7068.8568
$ws.Range("M122").Value = -178827.82
$ws.Range("N122").Value = -11968.8568
